$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare two brand-new rows (13, 14) with the same look as the
# existing data rows (border + number format) by copying row 12's format
# down before writing values into them. ---
$ws.Range("A12:D12").Copy()
$ws.Range("A13:D14").PasteSpecial(-4122)

# --- Rewrite the whole data block (rows 2-14) with the final values,
# which reorders the existing records and appends the two new ones. ---
$ws.Cells.Item(2,1).Value = "5612603000A"
$ws.Cells.Item(2,2).Value = "SIM-M23001"
$ws.Cells.Item(2,3).Value = 1520000
$ws.Cells.Item(2,4).Value = 12.666666666666666

$ws.Cells.Item(3,1).Value = "5612603100A"
$ws.Cells.Item(3,2).Value = "SIM-M23002"
$ws.Cells.Item(3,3).Value = 1295000
$ws.Cells.Item(3,4).Value = 10.791666666666666

$ws.Cells.Item(4,1).Value = "Z0011377A"
$ws.Cells.Item(4,2).Value = "SIM-M23019"
$ws.Cells.Item(4,3).Value = 2000000
$ws.Cells.Item(4,4).Value = 16.669999999999987

$ws.Cells.Item(5,1).Value = "Z0011378A"
$ws.Cells.Item(5,2).Value = "SIM-M23020"
$ws.Cells.Item(5,3).Value = 1965000
$ws.Cells.Item(5,4).Value = 16.379999999999995

$ws.Cells.Item(6,1).Value = "T907055A"
$ws.Cells.Item(6,2).Value = "SIM-M23021"
$ws.Cells.Item(6,3).Value = 2200000
$ws.Cells.Item(6,4).Value = 18.330000000000013

$ws.Cells.Item(7,1).Value = "T46515AB"
$ws.Cells.Item(7,2).Value = "SIM-M23023"
$ws.Cells.Item(7,3).Value = 1165000
$ws.Cells.Item(7,4).Value = 58.25

$ws.Cells.Item(8,1).Value = "T46497AA"
$ws.Cells.Item(8,2).Value = "SIM-M23022"
$ws.Cells.Item(8,3).Value = 1500000
$ws.Cells.Item(8,4).Value = 75

$ws.Cells.Item(9,1).Value = "5611012633A"
$ws.Cells.Item(9,2).Value = "SIM-M23026"
$ws.Cells.Item(9,3).Value = 1346000
$ws.Cells.Item(9,4).Value = 67.3

$ws.Cells.Item(10,1).Value = "5611019231A"
$ws.Cells.Item(10,2).Value = "SIM-M23027"
$ws.Cells.Item(10,3).Value = 1270000
$ws.Cells.Item(10,4).Value = 84.666666666666671

$ws.Cells.Item(11,1).Value = "5611019330A"
$ws.Cells.Item(11,2).Value = "SIM-M23028"
$ws.Cells.Item(11,3).Value = 1270000
$ws.Cells.Item(11,4).Value = 84.666666666666671

$ws.Cells.Item(12,1).Value = "T46511AB"
$ws.Cells.Item(12,2).Value = "SIM-M23048"
$ws.Cells.Item(12,3).Value = 1200000
$ws.Cells.Item(12,4).Value = 60

$ws.Cells.Item(13,1).Value = "Z0009680A"
$ws.Cells.Item(13,2).Value = "SIM-M23025"
$ws.Cells.Item(13,3).Value = 1800000
$ws.Cells.Item(13,4).Value = 120

$ws.Cells.Item(14,1).Value = "Z0009775A"
$ws.Cells.Item(14,2).Value = "SIM-M23052"
$ws.Cells.Item(14,3).Value = 1200000
$ws.Cells.Item(14,4).Value = 80

# --- A9 ("5611012633A") gets wrap text + vertically centered alignment,
# with its border cleared. ---
$ws.Range("A9").WrapText = $true
$ws.Range("A9").VerticalAlignment = -4108
$ws.Range("A9").Borders.LineStyle = -4142

# --- Column A widens slightly to fit the new content. ---
$ws.Columns.Item(1).ColumnWidth = 13.21875

# --- Selection ends up on A9. ---
$ws.Range("A9").Select()
